$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column D (Price) as Text so numeric-looking values are not
# auto-converted by Excel into numbers (the column holds strings like "30.191.57").
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '30.191.57'
$ws.Range("E2").Value = '  +8.15%  '

# Row 3
$ws.Range("D3").Value = '1.876.51'
$ws.Range("E3").Value = '  +5.43%  '

# Row 4
$ws.Range("D4").Value = '0.9978'
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$ws.Range("D5").Value = '250.00'
$ws.Range("E5").Value = '  +2.73%  '

# Row 6
$ws.Range("E6").Value = '  +0.17%  '

# Row 7
$ws.Range("D7").Value = '0.4985'
$ws.Range("E7").Value = '  +2.29%  '

# Row 8
$ws.Range("D8").Value = '0.2861'
$ws.Range("E8").Value = '  +7.84%  '

# Row 9
$ws.Range("D9").Value = '0.06597'
$ws.Range("E9").Value = '  +5.69%  '

# Row 10
$ws.Range("D10").Value = '1.867.59'
$ws.Range("E10").Value = '  +5.05%  '

# Row 11
$ws.Range("D11").Value = '17.13'
$ws.Range("E11").Value = '  +4.70%  '

# Row 12
$ws.Range("D12").Value = '0.07207'
$ws.Range("E12").Value = '  +3.39%  '

# Row 13
$ws.Range("D13").Value = '0.6644'
$ws.Range("E13").Value = '  +7.75%  '

# Row 14
$ws.Range("D14").Value = '85.40'
$ws.Range("E14").Value = '  +7.99%  '

# Row 15
$ws.Range("D15").Value = '4.831'
$ws.Range("E15").Value = '  +5.34%  '

# Row 16
$ws.Range("D16").Value = '30.154.26'
$ws.Range("E16").Value = '  +8.11%  '

# Row 17
$ws.Range("D17").Value = '0.9982'
$ws.Range("E17").Value = '  +0.14%  '

# Row 18
$ws.Range("D18").Value = '12.93'
$ws.Range("E18").Value = '  +9.96%  '

# Row 19
$ws.Range("D19").Value = '0.000007546'
$ws.Range("E19").Value = '  +4.89%  '

# Row 20
$ws.Range("D20").Value = '0.9973'
$ws.Range("E20").Value = '  +0.12%  '

# Row 21
$ws.Range("D21").Value = '2.106.28'
$ws.Range("E21").Value = '  +5.65%  '

# Row 22
$ws.Range("D22").Value = '4.787'
$ws.Range("E22").Value = '  +4.57%  '

# Row 23
$ws.Range("D23").Value = '9.082'
$ws.Range("E23").Value = '  +5.49%  '

# Row 24
$ws.Range("D24").Value = '5.519'
$ws.Range("E24").Value = '  +6.65%  '

# Row 25
$ws.Range("B25").Value = 'BitcoinCash'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D25").Value = '138.19'
$ws.Range("E25").Value = '  +25.72%  '

# Row 26
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '145.18'
$ws.Range("E26").Value = '  +2.68%  '

# Row 27
$ws.Range("D27").Value = '16.78'
$ws.Range("E27").Value = '  +7.67%  '

# Row 28
$ws.Range("D28").Value = '1.960'
$ws.Range("E28").Value = '  +4.41%  '

# Row 29
$ws.Range("D29").Value = '1.392'
$ws.Range("E29").Value = '  +0.72%  '

# Row 30
$ws.Range("D30").Value = '4.265'
$ws.Range("E30").Value = '  +6.19%  '

# Row 31
$ws.Range("D31").Value = '0.08642'
$ws.Range("E31").Value = '  +4.62%  '

# Row 32
$ws.Range("D32").Value = '3.912'
$ws.Range("E32").Value = '  +4.17%  '

# Row 33
$ws.Range("D33").Value = '0.05100'
$ws.Range("E33").Value = '  +8.50%  '

# Row 34
$ws.Range("D34").Value = '1.135'
$ws.Range("E34").Value = '  +8.52%  '

# Row 35
$ws.Range("D35").Value = '0.6930'
$ws.Range("E35").Value = '  +9.67%  '

# Row 36
$ws.Range("D36").Value = '2.691'
$ws.Range("E36").Value = '  +2.68%  '

# Row 37
$ws.Range("D37").Value = '2.334'
$ws.Range("E37").Value = '  +13.91%  '

# Row 38
$ws.Range("D38").Value = '2.754'
$ws.Range("E38").Value = '  +6.17%  '

# Row 39
$ws.Range("D39").Value = '0.9625'
$ws.Range("E39").Value = '  +2.38%  '

# Row 40
$ws.Range("D40").Value = '0.01638'
$ws.Range("E40").Value = '  +7.26%  '

# Row 41
$ws.Range("D41").Value = '6.140'
$ws.Range("E41").Value = '  +5.99%  '

# Row 42
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '104.14'
$ws.Range("E42").Value = '  +4.34%  '

# Row 43
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = '0.9979'
$ws.Range("E43").Value = '  +0.10%  '

# Row 44
$ws.Range("D44").Value = '0.4212'
$ws.Range("E44").Value = '  +7.45%  '

# Row 45
$ws.Range("D45").Value = '7.465'
$ws.Range("E45").Value = '  +6.77%  '

# Row 46
$ws.Range("D46").Value = '0.1259'
$ws.Range("E46").Value = '  +6.17%  '

# Row 47
$ws.Range("D47").Value = '0.05638'
$ws.Range("E47").Value = '  +4.30%  '

# Row 48
$ws.Range("D48").Value = '32.67'
$ws.Range("E48").Value = '  +7.92%  '

# Row 49
$ws.Range("D49").Value = '8.312'
$ws.Range("E49").Value = '  +5.14%  '

# Row 50
$ws.Range("D50").Value = '0.3742'
$ws.Range("E50").Value = '  +8.55%  '

# Row 51
$ws.Range("D51").Value = '1.345'
$ws.Range("E51").Value = '  +6.00%  '
